# Crowdin sync: Japanese (ja) localization sheet.
# The English (column E) reference translations that used to sit alongside
# the Japanese strings in rows 26-44 are dropped from this locale sheet
# (Crowdin only keeps the JA values here), except the very last entry
# (row 46, "AllDismissItem") which keeps its English value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToClearColumnE = @(26,27,28,29,30,31,33,34,35,36,38,39,40,41,43,44)

foreach ($r in $rowsToClearColumnE) {
    $ws.Cells.Item($r, 5).ClearContents()
}

# Match the saved selection state (cursor left on E47 after the edit).
$ws.Range("E47").Select() | Out-Null
